$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.118.84"
$ws.Range("E2").Value = "  +5.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.69"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("E4").Value = "  -0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.78"
$ws.Range("E5").Value = "  +4.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5207"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4080"
$ws.Range("E8").Value = "  +4.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08561"
$ws.Range("E9").Value = "  +2.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.02"
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.126"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.35"
$ws.Range("E12").Value = "  +9.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.423"
$ws.Range("E13").Value = "  +3.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.921.97"
$ws.Range("E14").Value = "  +2.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.410"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.34"
$ws.Range("E17").Value = "  +4.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06688"
$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.015"
$ws.Range("E22").Value = "  +1.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.134.12"
$ws.Range("E23").Value = "  +5.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.202"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.144.71"
$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.10"
$ws.Range("E27").Value = "  +2.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.73"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.443"
$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.04"
$ws.Range("E30").Value = "  +1.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.082"
$ws.Range("E31").Value = "  +3.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1062"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.058"
$ws.Range("E33").Value = "  +5.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.637"
$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02491"
$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.232"
$ws.Range("E38").Value = "  +4.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.176"
$ws.Range("E39").Value = "  +2.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.914"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6550"
$ws.Range("E41").Value = "  +2.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.249"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.63"
$ws.Range("E43").Value = "  +4.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6153"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.22"
$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("E46").Value = "  +2.20%  "

$ws.Range("E47").Value = "  +3.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.249"
$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.35"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.173"
$ws.Range("E50").Value = "  +10.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.76"
$ws.Range("E51").Value = "  +4.30%  "

